$d = $word.ActiveDocument
$d.TrackRevisions = $false

# --- Step 1: insert the new protected-category list before "ANY COMPUTER PROGRAMMER" ---
$anchor = $d.Content
$found = $anchor.Find.Execute("ANY COMPUTER PROGRAMMER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Anchor text not found: ANY COMPUTER PROGRAMMER" }
$insertStart = $anchor.Start
$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertBefore("ANY ROOFER XOR ANY REAL ESTATE AGENT XOR ANY PLUMBER XOR ANY WORKER XOR           ANY UNEMPLOYED PERSON XOR ANY CHILD XOR ANY PARENT XOR ANY SINGLE PERSON XOR          ")

# --- Step 2: colour the newly inserted text segment by segment ---
$seg = $d.Range($insertStart + 0, $insertStart + 10)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 10, $insertStart + 14)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 15, $insertStart + 37)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 37, $insertStart + 40)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 41, $insertStart + 53)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 53, $insertStart + 56)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 57, $insertStart + 67)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 67, $insertStart + 71)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 82, $insertStart + 103)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 103, $insertStart + 107)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 108, $insertStart + 117)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 117, $insertStart + 121)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 122, $insertStart + 132)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 132, $insertStart + 136)
$seg.Font.Color = 15773696
$seg = $d.Range($insertStart + 137, $insertStart + 154)
$seg.Font.Color = 255
$seg = $d.Range($insertStart + 154, $insertStart + 158)
$seg.Font.Color = 15773696

# --- Step 3: merge the "ANY " + "PROTECTED INTERNATIONAL CRIMINAL COURT WITNESS" runs ---
$d.Content.Find.Execute("ANY PROTECTED INTERNATIONAL CRIMINAL COURT WITNESS", $true, $false, $false, $false, $false, $true, 1, $false, "ANY PROTECTED INTERNATIONAL CRIMINAL COURT WITNESS", 2) | Out-Null

# --- Step 4: merge the "ANY" + " INTERPOL PROTECTED WITNESS" runs ---
$d.Content.Find.Execute("ANY INTERPOL PROTECTED WITNESS", $true, $false, $false, $false, $false, $true, 1, $false, "ANY INTERPOL PROTECTED WITNESS", 2) | Out-Null

Write-Output "done"
